$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.311647653579712
$ws.Range("B1").Value = 2.47678804397583
$ws.Range("C1").Value = 2.084283351898193
$ws.Range("D1").Value = 2.167176246643066
$ws.Range("E1").Value = 2.506239414215088
